$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.189.39'
$ws.Range('E2').Value = '  +0.67%  '

$ws.Range('D3').Value = '3.947.07'
$ws.Range('E3').Value = '  +4.25%  '

$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').Value = '''471.50'
$ws.Range('E5').Value = '  +8.98%  '

$ws.Range('D6').Value = '''146.06'
$ws.Range('E6').Value = '  +4.19%  '

$ws.Range('E7').Value = '  +0.77%  '

$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('E9').Value = '  -0.17%  '

$ws.Range('D10').Value = '''0.163'
$ws.Range('E10').Value = '  +6.63%  '

$ws.Range('D11').Value = '''0.0000336'
$ws.Range('E11').Value = '  +6.26%  '

$ws.Range('D12').Value = '''43.32'
$ws.Range('E12').Value = '  +0.95%  '

$ws.Range('D13').Value = '4.569.83'
$ws.Range('E13').Value = '  +4.57%  '

$ws.Range('D14').Value = '''10.35'
$ws.Range('E14').Value = '  -0.72%  '

$ws.Range('D15').Value = '''15.29'
$ws.Range('E15').Value = '  +1.94%  '

$ws.Range('D16').Value = '3.931.28'
$ws.Range('E16').Value = '  +2.63%  '

$ws.Range('E17').Value = '  -0.12%  '

$ws.Range('D18').Value = '''19.86'
$ws.Range('E18').Value = '  -0.40%  '

$ws.Range('E19').Value = '  +2.11%  '

$ws.Range('D20').Value = '67.476.46'
$ws.Range('E20').Value = '  +1.24%  '

$ws.Range('D21').Value = '''438.89'
$ws.Range('E21').Value = '  +7.00%  '

$ws.Range('D22').Value = '''3.41'
$ws.Range('E22').Value = '  +4.62%  '

$ws.Range('D23').Value = '''14.55'
$ws.Range('E23').Value = '  -1.69%  '

$ws.Range('D24').Value = '''87.71'
$ws.Range('E24').Value = '  +2.63%  '

$ws.Range('D25').Value = '''3.64'
$ws.Range('E25').Value = '  +8.62%  '

$ws.Range('E26').Value = '  +5.93%  '

$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').Value = '''10.17'
$ws.Range('E27').Value = '  +3.75%  '

$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '''9.75'
$ws.Range('E28').Value = '  +0.88%  '

$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = '''723.29'
$ws.Range('E29').Value = '  +1.47%  '

$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '''13.55'
$ws.Range('E30').Value = '  -2.08%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.132'
$ws.Range('E31').Value = '  -2.05%  '

$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '''2.84'
$ws.Range('E32').Value = '  +2.94%  '

$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = '''42.86'
$ws.Range('E33').Value = '  +2.71%  '

$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '''0.153'
$ws.Range('E34').Value = '  +0.49%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''57.91'
$ws.Range('E35').Value = '  +3.23%  '

$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = '''0.999'
$ws.Range('E36').Value = '  +0.02%  '

$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0785'
$ws.Range('E37').Value = '  +15.47%  '

$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '''5.36'
$ws.Range('E38').Value = '  -5.72%  '

$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0478'
$ws.Range('E39').Value = '  +0.35%  '

$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').Value = '''3.06'
$ws.Range('E40').Value = '  +4.12%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '''0.142'
$ws.Range('E41').Value = '  +0.54%  '

$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  +0.11%  '

$ws.Range('E43').Value = '  +3.69%  '

$ws.Range('B44').Value = 'LidoDAOToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D44').Value = '''3.49'
$ws.Range('E44').Value = '  +4.54%  '

$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = '''2.56'
$ws.Range('E45').Value = '  -8.81%  '

$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '''2.20'
$ws.Range('E46').Value = '  +5.22%  '

$ws.Range('D47').Value = '''2.79'
$ws.Range('E47').Value = '  +3.62%  '

$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''146.50'
$ws.Range('E48').Value = '  +2.73%  '

$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '''3.14'
$ws.Range('E49').Value = '  -4.60%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '''2.86'
$ws.Range('E50').Value = '  +1.49%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''25.91'
$ws.Range('E51').Value = '  +3.31%  '
